# Generate Report for Handback
# The 211f772b-...md file has been handed back (in sync with en-US) for both
# the zh-cn and de-de languages. This updates the Overview sheet status, the
# per-language sheets' Status / Latest Target File / Latest Handback File /
# Latest Handback DateTime columns, and adds the corresponding hyperlinks.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: file 211f772b is now handed back (both languages)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2 corresponds to file 211f772b-b4a1-4978-9a14-73fdbae09532.md
$wsZh.Range("C2").Value = $statusHandedBack

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d6ca24078cba4337427b9bb0c1ff7e6a770be1f/e2e/211f772b-b4a1-4978-9a14-73fdbae09532.md", "", "", "211f772b-b4a1-4978-9a14-73fdbae09532.md") | Out-Null
$wsZh.Range("F2").Font.Name = "Calibri"
$wsZh.Range("F2").Font.Size = 11
$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("F2").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6177faad2eac19c89bd54677de541fe61767987c/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/hb/211f772b-b4a1-4978-9a14-73fdbae09532.d2b3ef8249b7c323e54573d278d309a551062cf8.zh-cn.xlf", "", "", "211f772b-b4a1-4978-9a14-73fdbae09532.d2b3ef8249b7c323e54573d278d309a551062cf8.zh-cn.xlf") | Out-Null
$wsZh.Range("G2").Font.Name = "Calibri"
$wsZh.Range("G2").Font.Size = 11
$wsZh.Range("G2").Font.Underline = 2
$wsZh.Range("G2").Font.Color = 15570276

$wsZh.Range("H2").Value = "2016-03-17 20:29:05"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2 corresponds to file 211f772b-b4a1-4978-9a14-73fdbae09532.md
$wsDe.Range("C2").Value = $statusHandedBack

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/1d6ca24078cba4337427b9bb0c1ff7e6a770be1f/e2e/211f772b-b4a1-4978-9a14-73fdbae09532.md", "", "", "211f772b-b4a1-4978-9a14-73fdbae09532.md") | Out-Null
$wsDe.Range("F2").Font.Name = "Calibri"
$wsDe.Range("F2").Font.Size = 11
$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("F2").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/8f4bdd23f1707fbba67b63d55cf075ccaf9a0585/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/hb/211f772b-b4a1-4978-9a14-73fdbae09532.d2b3ef8249b7c323e54573d278d309a551062cf8.de-de.xlf", "", "", "211f772b-b4a1-4978-9a14-73fdbae09532.d2b3ef8249b7c323e54573d278d309a551062cf8.de-de.xlf") | Out-Null
$wsDe.Range("G2").Font.Name = "Calibri"
$wsDe.Range("G2").Font.Size = 11
$wsDe.Range("G2").Font.Underline = 2
$wsDe.Range("G2").Font.Color = 15570276

$wsDe.Range("H2").Value = "2016-03-17 20:29:11"
